$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 11-13 (Neutrophils sending-cluster block no longer present)
$ws.Rows.Item(11).Resize(3).Delete()

# Update rows 2-10 with refreshed TPM-derived values and cluster labels
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ngfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.1450045
$ws.Range("H2").Value = 14.290009
$ws.Range("I2").Value = 0.8119737125238713
$ws.Range("J2").Value = 0.7990590344890214
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.337665
$ws.Range("N2").Value = 0.67533
$ws.Range("O2").Value = 0.01438349055657064
$ws.Range("P2").Value = 0.0143657952272707
$ws.Range("Q2").Value = 2.4126179444925
$ws.Range("R2").Value = 9.650471777969999
$ws.Range("S2").Value = 0.0116790162262707
$ws.Range("T2").Value = 0.01147911846396992

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ngfr"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.1450045
$ws.Range("H3").Value = 14.290009
$ws.Range("I3").Value = 0.8119737125238713
$ws.Range("J3").Value = 0.7990590344890214
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05783366666666667
$ws.Range("N3").Value = 0.173501
$ws.Range("O3").Value = 0.002463536340313192
$ws.Range("P3").Value = 0.003690758351808291
$ws.Range("Q3").Value = 0.4132218085848333
$ws.Range("R3").Value = 2.479330851509
$ws.Range("S3").Value = 0.002000326748181574
$ws.Range("T3").Value = 0.002949133805128225

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ngfr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.1450045
$ws.Range("H4").Value = 14.290009
$ws.Range("I4").Value = 0.8119737125238713
$ws.Range("J4").Value = 0.7990590344890214
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.0803745
$ws.Range("N4").Value = 46.160749
$ws.Range("O4").Value = 0.9831529731031161
$ws.Range("P4").Value = 0.981943446420921
$ws.Range("Q4").Value = 164.9093796641852
$ws.Range("R4").Value = 659.6375186567409
$ws.Range("S4").Value = 0.798294369549419
$ws.Range("T4").Value = 0.7846307822199232

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ngfr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.284443
$ws.Range("H5").Value = 0.853329
$ws.Range("I5").Value = 0.0323247156403369
$ws.Range("J5").Value = 0.04771587245616726
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.337665
$ws.Range("N5").Value = 0.67533
$ws.Range("O5").Value = 0.01438349055657064
$ws.Range("P5").Value = 0.0143657952272707
$ws.Range("Q5").Value = 0.096046445595
$ws.Range("R5").Value = 0.57627867357
$ws.Range("S5").Value = 0.0004649422421566169
$ws.Range("T5").Value = 0.000685476452795865

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ngfr"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.284443
$ws.Range("H6").Value = 0.853329
$ws.Range("I6").Value = 0.0323247156403369
$ws.Range("J6").Value = 0.04771587245616726
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05783366666666667
$ws.Range("N6").Value = 0.173501
$ws.Range("O6").Value = 0.002463536340313192
$ws.Range("P6").Value = 0.003690758351808291
$ws.Range("Q6").Value = 0.01645038164766667
$ws.Range("R6").Value = 0.148053434829
$ws.Range("S6").Value = 0.00007963311167026017
$ws.Range("T6").Value = 0.0001761077547814185

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ngfr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.284443
$ws.Range("H7").Value = 0.853329
$ws.Range("I7").Value = 0.0323247156403369
$ws.Range("J7").Value = 0.04771587245616726
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.0803745
$ws.Range("N7").Value = 46.160749
$ws.Range("O7").Value = 0.9831529731031161
$ws.Range("P7").Value = 0.981943446420921
$ws.Range("Q7").Value = 6.5650509639035
$ws.Range("R7").Value = 39.390305783421
$ws.Range("S7").Value = 0.03178014028651002
$ws.Range("T7").Value = 0.04685428824858998

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ntf3"
$ws.Range("C8").Value = "Ngfr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.370104
$ws.Range("H8").Value = 2.740208
$ws.Range("I8").Value = 0.1557015718357919
$ws.Range("J8").Value = 0.1532250930548114
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.337665
$ws.Range("N8").Value = 0.67533
$ws.Range("O8").Value = 0.01438349055657064
$ws.Range("P8").Value = 0.0143657952272707
$ws.Range("Q8").Value = 0.46263616716
$ws.Range("R8").Value = 1.85054466864
$ws.Range("S8").Value = 0.002239532088143317
$ws.Range("T8").Value = 0.002201200310504918

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ntf3"
$ws.Range("C9").Value = "Ngfr"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.370104
$ws.Range("H9").Value = 2.740208
$ws.Range("I9").Value = 0.1557015718357919
$ws.Range("J9").Value = 0.1532250930548114
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05783366666666667
$ws.Range("N9").Value = 0.173501
$ws.Range("O9").Value = 0.002463536340313192
$ws.Range("P9").Value = 0.003690758351808291
$ws.Range("Q9").Value = 0.07923813803466666
$ws.Range("R9").Value = 0.475428828208
$ws.Range("S9").Value = 0.0003835764804613583
$ws.Range("T9").Value = 0.0005655167918986477

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ntf3"
$ws.Range("C10").Value = "Ngfr"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.370104
$ws.Range("H10").Value = 2.740208
$ws.Range("I10").Value = 0.1557015718357919
$ws.Range("J10").Value = 0.1532250930548114
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 23.0803745
$ws.Range("N10").Value = 46.160749
$ws.Range("O10").Value = 0.9831529731031161
$ws.Range("P10").Value = 0.981943446420921
$ws.Range("Q10").Value = 31.622513423948
$ws.Range("R10").Value = 126.490053695792
$ws.Range("S10").Value = 0.1530784632671872
$ws.Range("T10").Value = 0.1504583759524078

